# New changes - 6/30/2020
#
# Settings sheet:
#   - Removes the "OrchestratorQueueName / ProcessABCQueue" row.
#   - "logF_BusinessProcessName" row moves up to row 2 (keeps its ht=30 / wrapped description).
#   - "MappingFilePath" row moves up to row 3, and its value is repointed from the
#     old dev machine path to the new one.
#   - Three new named settings are appended (rows 4-6): MappingFilePath_BankReport_SheetName,
#     MappingFilePath_RemittanceReport_SheetName and PNC_BankReport_SheetName.
#
# Constants sheet:
#   - The two blank spacer rows are removed so the constants are contiguous.
#   - The MaxRetryNumber value cell is left-aligned.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Delete the "OrchestratorQueueName" row (row 2) outright - its content does not
# survive the edit anywhere else in the sheet.
$settings.Rows.Item(2).Delete()

# The sheet now reads: row2=logF_BusinessProcessName(ht=30), row3=blank(was4),
# row4=blank(was5), row5=blank(was6), row6=MappingFilePath(was7). Collapse the
# blank rows so logF stays at row 2 and MappingFilePath lands on row 3.
$settings.Rows.Item(3).Delete()
$settings.Rows.Item(3).Delete()

# Update the MappingFilePath value to the new workstation path.
$settings.Cells.Item(3, 1).Value = "MappingFilePath"
$settings.Cells.Item(3, 2).Value = "C:\Users\Hp\Documents\UiPath\AR2.0\Data\MappingSheet.xlsx"

# Append the three new settings rows.
$settings.Cells.Item(4, 1).Value = "MappingFilePath_BankReport_SheetName"
$settings.Cells.Item(4, 2).Value = "Bank_Report"

$settings.Cells.Item(5, 1).Value = "MappingFilePath_RemittanceReport_SheetName"
$settings.Cells.Item(5, 2).Value = "Remittance_Report"

$settings.Cells.Item(6, 1).Value = "PNC_BankReport_SheetName"
$settings.Cells.Item(6, 2).Value = "Bank"

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# Remove the two blank spacer rows (originally row 5, then row 3) so every
# constant lives in a contiguous block starting at row 2.
$constants.Rows.Item(5).Delete()
$constants.Rows.Item(3).Delete()

# Left-align the MaxRetryNumber value cell.
$constants.Cells.Item(2, 2).HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Restore selections / active sheet to match the saved view state.
# ---------------------------------------------------------------------------
$constants.Range("B2").Select()
$settings.Activate()
$settings.Range("A6").Select()
